# Finalized Experiments with Participant Generation
# Rename worksheets (new participant-generation timestamps) and update
# the generated stimulus-file names / condition order values that were
# regenerated along with them.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (order preserved) ---------------------------------
$wsGNG  = $wb.Worksheets.Item(1)
$wsNB   = $wb.Worksheets.Item(2)
$wsRS   = $wb.Worksheets.Item(3)
$wsTOL  = $wb.Worksheets.Item(4)
$wsvSAT = $wb.Worksheets.Item(5)

$wsGNG.Name  = "GNG_TO-16502912950093892"
$wsNB.Name   = "NB_TO-1650291297313209"
$wsRS.Name   = "RS_TO-16502912973142076"
$wsTOL.Name  = "TOL_TO-16502912973773"
$wsvSAT.Name = "vSAT_TO-16502912974481826"

# --- Sheet 1: GNG ------------------------------------------------------
$wsGNG.Range("B2").Value = "go_stims-16502912949601977.csv"
$wsGNG.Range("B3").Value = "GNG_stims-16502912949767191.csv"
$wsGNG.Range("B4").Value = "go_stims-16502912949787767.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16502912950083568.csv"

# --- Sheet 2: NB ---------------------------------------------------------
$wsNB.Range("B2").Value = "OB-16502912961915188.csv"
$wsNB.Range("B3").Value = "OB-16502912963880095.csv"
$wsNB.Range("B4").Value = "ZB-match_7-16502912951859667.csv"
$wsNB.Range("B5").Value = "TB-16502912968046935.csv"
$wsNB.Range("B6").Value = "ZB-match_1-16502912955048096.csv"
$wsNB.Range("B7").Value = "OB-1650291296337462.csv"
$wsNB.Range("B8").Value = "TB-16502912972597625.csv"
$wsNB.Range("B9").Value = "TB-16502912972948935.csv"
$wsNB.Range("B10").Value = "ZB-match_5-16502912956519852.csv"

# --- Sheet 3: RS ---------------------------------------------------------
$wsRS.Range("B2").Value = "eyes open"
$wsRS.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---------------------------------------------------------
$wsTOL.Range("B2").Value = "MM_stims-16502912973287969.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16502912973162494.csv"
$wsTOL.Range("B4").Value = "MM_stims-1650291297359939.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16502912973298056.csv"
$wsTOL.Range("B6").Value = "MM_stims-16502912973763.csv"
$wsTOL.Range("B7").Value = "ZM_stims-1650291297359939.csv"

# --- Sheet 5: vSAT ---------------------------------------------------------
$wsvSAT.Range("B2").Value = "SAT_stims-1650291297381412.csv"
$wsvSAT.Range("B3").Value = "vSAT_stims-1650291297417404.csv"
$wsvSAT.Range("B4").Value = "SAT_stims-16502912973917892.csv"
$wsvSAT.Range("B5").Value = "vSAT_stims-16502912974327638.csv"
